$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text, even when it looks like
# a pure number (e.g. "0.997", "1.00") or would otherwise be mis-parsed.
# We force the cell to a text number-format before assigning the value so
# Excel's type auto-detection doesn't coerce it into a numeric cell, then
# restore the original "Normal" cell style so no formatting residue is left
# behind on the cell.
function Set-CellNumericText {
    param($sheet, [string]$addr, [string]$value)
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Helper: write a plain text value (not at risk of numeric auto-detection).
function Set-CellText {
    param($sheet, [string]$addr, [string]$value)
    $sheet.Range($addr).Value = $value
}

Set-CellNumericText $ws "D2" "57.716.08"
Set-CellText $ws "E2" "  -0.48%  "
Set-CellNumericText $ws "D3" "2.438.69"
Set-CellText $ws "E3" "  -1.11%  "
Set-CellText $ws "E4" "  +0.00%  "
Set-CellNumericText $ws "D5" "508.34"
Set-CellText $ws "E5" "  -1.71%  "
Set-CellNumericText $ws "D6" "129.25"
Set-CellText $ws "E6" "  -2.00%  "
Set-CellNumericText $ws "D7" "0.997"
Set-CellText $ws "E7" "  -0.30%  "
Set-CellNumericText $ws "D8" "0.550"
Set-CellText $ws "E8" "  -1.30%  "
Set-CellNumericText $ws "D9" "2.456.56"
Set-CellText $ws "E9" "  -0.44%  "
Set-CellText $ws "E10" "  -0.14%  "
Set-CellNumericText $ws "D11" "0.0954"
Set-CellText $ws "E11" "  -3.95%  "
Set-CellNumericText $ws "D12" "5.17"
Set-CellText $ws "E12" "  -4.11%  "
Set-CellNumericText $ws "D13" "0.330"
Set-CellText $ws "E13" "  -3.21%  "
Set-CellNumericText $ws "D14" "2.869.08"
Set-CellText $ws "E14" "  -1.21%  "
Set-CellNumericText $ws "D15" "57.640.47"
Set-CellText $ws "E15" "  -0.46%  "
Set-CellNumericText $ws "D16" "21.91"
Set-CellText $ws "E16" "  -0.59%  "
Set-CellNumericText $ws "D17" "0.0000133"
Set-CellText $ws "E17" "  -2.70%  "
Set-CellNumericText $ws "D18" "2.438.65"
Set-CellText $ws "E18" "  -1.63%  "
Set-CellNumericText $ws "D19" "10.47"
Set-CellText $ws "E19" "  -3.52%  "
Set-CellNumericText $ws "D20" "4.12"
Set-CellText $ws "E20" "  -1.17%  "
Set-CellNumericText $ws "D21" "315.18"
Set-CellText $ws "E21" "  -1.19%  "
Set-CellNumericText $ws "D23" "5.68"
Set-CellText $ws "E23" "  -1.47%  "
Set-CellNumericText $ws "D24" "63.38"
Set-CellText $ws "E24" "  -1.66%  "
Set-CellNumericText $ws "D25" "0.406"
Set-CellText $ws "E25" "  -0.66%  "
Set-CellText $ws "E26" "  -0.37%  "
Set-CellNumericText $ws "D27" "0.160"
Set-CellText $ws "E27" "  -1.16%  "
Set-CellNumericText $ws "D28" "7.27"
Set-CellText $ws "E28" "  -1.74%  "
Set-CellNumericText $ws "D29" "169.72"
Set-CellText $ws "E29" "  +3.05%  "
Set-CellNumericText $ws "D30" "0.0₃0724"
Set-CellText $ws "E30" "  -3.02%  "
Set-CellNumericText $ws "D31" "6.26"
Set-CellText $ws "E31" "  -1.86%  "
Set-CellText $ws "B32" "Fetch.AI"
Set-CellText $ws "C32" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellNumericText $ws "D32" "1.16"
Set-CellText $ws "E32" "  +2.39%  "
Set-CellText $ws "B33" "PancakeSwap"
Set-CellText $ws "C33" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellNumericText $ws "D33" "1.66"
Set-CellText $ws "E33" "  -2.61%  "
Set-CellNumericText $ws "D34" "0.997"
Set-CellText $ws "E34" "  -0.06%  "
Set-CellNumericText $ws "D35" "0.997"
Set-CellText $ws "E35" "  -0.24%  "
Set-CellNumericText $ws "D36" "17.71"
Set-CellText $ws "E36" "  -2.40%  "
Set-CellNumericText $ws "D37" "1.27"
Set-CellText $ws "E37" "  -4.66%  "
Set-CellNumericText $ws "D38" "3.94"
Set-CellText $ws "E38" "  -0.48%  "
Set-CellNumericText $ws "D39" "36.28"
Set-CellText $ws "E39" "  -0.56%  "
Set-CellNumericText $ws "D40" "1.46"
Set-CellText $ws "E40" "  -1.70%  "
Set-CellNumericText $ws "D41" "0.767"
Set-CellText $ws "E41" "  -3.18%  "
Set-CellNumericText $ws "D42" "272.78"
Set-CellText $ws "E42" "  -0.64%  "
Set-CellText $ws "B43" "RenderToken"
Set-CellText $ws "C43" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellNumericText $ws "D43" "5.01"
Set-CellText $ws "E43" "  +0.94%  "
Set-CellText $ws "B44" "Filecoin"
Set-CellText $ws "C44" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellNumericText $ws "D44" "3.39"
Set-CellText $ws "E44" "  -2.78%  "
Set-CellNumericText $ws "D45" "0.582"
Set-CellText $ws "E45" "  -1.32%  "
Set-CellNumericText $ws "D46" "0.0909"
Set-CellText $ws "E46" "  -0.11%  "
Set-CellNumericText $ws "D47" "120.26"
Set-CellText $ws "E47" "  -5.29%  "
Set-CellNumericText $ws "D48" "0.0485"
Set-CellText $ws "E48" "  -1.20%  "
Set-CellNumericText $ws "D49" "17.21"
Set-CellText $ws "E49" "  -2.97%  "
Set-CellNumericText $ws "D50" "0.0210"
Set-CellText $ws "E50" "  -2.33%  "
Set-CellNumericText $ws "D51" "16.66"
Set-CellText $ws "E51" "  -2.39%  "

Write-Host "Updated cryptos list"
